# This script updates column B ("pred") values on Sheet1 for a set of rows,
# per the "keep highest val strategy" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row=11;  Old="AAA"; New="A"},
    @{Row=19;  Old="BBB"; New="A"},
    @{Row=30;  Old="AA";  New="A"},
    @{Row=38;  Old="A";   New="BBB"},
    @{Row=46;  Old="BB";  New="B"},
    @{Row=58;  Old="BBB"; New="BB"},
    @{Row=59;  Old="BB";  New="B"},
    @{Row=69;  Old="B";   New="BB"},
    @{Row=70;  Old="B";   New="BB"},
    @{Row=74;  Old="BB";  New="B"},
    @{Row=80;  Old="BB";  New="B"},
    @{Row=81;  Old="AA";  New="AAA"},
    @{Row=87;  Old="A";   New="BBB"},
    @{Row=94;  Old="BBB"; New="A"},
    @{Row=99;  Old="BBB"; New="A"},
    @{Row=116; Old="A";   New="BBB"},
    @{Row=124; Old="A";   New="BBB"},
    @{Row=136; Old="BB";  New="BBB"},
    @{Row=143; Old="AA";  New="BBB"},
    @{Row=146; Old="BBB"; New="BB"},
    @{Row=159; Old="AA";  New="A"},
    @{Row=160; Old="AA";  New="A"},
    @{Row=166; Old="BBB"; New="A"},
    @{Row=167; Old="BBB"; New="A"},
    @{Row=201; Old="BB";  New="BBB"},
    @{Row=204; Old="AAA"; New="A"},
    @{Row=207; Old="A";   New="BBB"},
    @{Row=233; Old="BBB"; New="BB"},
    @{Row=281; Old="AA";  New="A"},
    @{Row=286; Old="BBB"; New="A"},
    @{Row=287; Old="BBB"; New="A"},
    @{Row=291; Old="A";   New="BBB"},
    @{Row=292; Old="A";   New="BBB"},
    @{Row=293; Old="A";   New="BBB"},
    @{Row=302; Old="A";   New="AA"},
    @{Row=306; Old="BBB"; New="A"},
    @{Row=307; Old="A";   New="BBB"},
    @{Row=308; Old="A";   New="BBB"},
    @{Row=315; Old="BBB"; New="BB"},
    @{Row=317; Old="B";   New="BB"},
    @{Row=323; Old="BB";  New="B"},
    @{Row=330; Old="BBB"; New="A"},
    @{Row=342; Old="BB";  New="BBB"},
    @{Row=348; Old="BB";  New="BBB"},
    @{Row=349; Old="B";   New="BB"},
    @{Row=355; Old="BBB"; New="BB"},
    @{Row=362; Old="BB";  New="B"},
    @{Row=375; Old="AA";  New="A"},
    @{Row=383; Old="BB";  New="B"},
    @{Row=387; Old="BB";  New="BBB"},
    @{Row=393; Old="BBB"; New="A"},
    @{Row=401; Old="BBB"; New="BB"}
)

foreach ($change in $changes) {
    $cell = $ws.Cells.Item($change.Row, 2)
    $cell.Value = $change.New
}
